$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (2-7) got reshuffled into a new order while the
# header (row 1) and the final "RHYS ap GRUFFYDD" row (row 8) stay put.
# Re-write columns A (title), B (timestamp), C (historical distance) and
# E (uri) for rows 2-7 with their new values; column D (time bucket) is the
# same "outside bucket range" value for all of these rows so it is left
# untouched.

$ws.Range("A2").Value = "Novel inspired by first eisteddfod and historic castle"
$ws.Range("B2").Value = "2018-03-17T00:00:00UTC"
$ws.Range("C2").Value = 307610
$ws.Range("E2").Value = "http://www.cambrian-news.co.uk/article.cfm?id=119720&headline=Novel%20inspired%20by%20first%20eisteddfod%20and%20historic%20castle&sectionIs=news&searchyear=2018"

$ws.Range("A3").Value = "Medieval legacy is there to be enjoyed"
$ws.Range("B3").Value = "2011-03-24T00:00:00UTC"
$ws.Range("C3").Value = 305060
$ws.Range("E3").Value = "https://www.walesonline.co.uk/lifestyle/showbiz/medieval-legacy-enjoyed-1848163"

$ws.Range("A4").Value = "The first eisteddfod - Christmas 1176"
$ws.Range("B4").Value = "2010-12-22T13:07:23UTC"
$ws.Range("C4").Value = 304968
$ws.Range("E4").Value = "https://www.bbc.co.uk/blogs/wales/entries/a1b7c602-5c9b-3bc1-ba13-fd64dc068ffc"

$ws.Range("A5").Value = "Nine foot replica of original National Eisteddfod chair is installed at the sitwe of Wales' very first cultural festival"
$ws.Range("B5").Value = "2015-03-18T16:52:24UTC"
$ws.Range("C5").Value = 306515
$ws.Range("E5").Value = "https://www.walesonline.co.uk/news/wales-news/nine-foot-replica-original-national-8870165"

$ws.Range("A6").Value = "Eisteddfod may have sprung from French idea"
$ws.Range("B6").Value = "2007-02-09T00:00:00UTC"
$ws.Range("C6").Value = 303556
$ws.Range("E6").Value = "https://www.walesonline.co.uk/news/wales-news/eisteddfod-sprung-french-idea-2273397"

$ws.Range("A7").Value = "Staff at Cardigan Castle's 1176 restaurant round off hectic year with busy Christmas period"
$ws.Range("B7").Value = "2016-12-19T09:51:13UTC"
$ws.Range("C7").Value = 307157
$ws.Range("E7").Value = "https://www.tivysideadvertiser.co.uk/news/14974769.restaurant-staff-set-to-turn-back-clock-with-christmas-servings/"
